# Update "Phi coc" (deposit fee) message and remove the second book row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update deposit fee text from "1112000 đồng" to "240000 đồng"
$ws.Range("B9").Value = "240000 đồng"

# Clear the second book row (title, author, call number) - row 13,
# copying row 14's (already-blank) formatting onto B13 so the border /
# alignment matches the other blank rows instead of keeping the old
# "header-like" centered/wrapped style.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13:D13").ClearContents()
